$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 109.9114832445916
$ws.Range("D2").Value = 3993.344853322108
$ws.Range("E2").Value = 1.354078223128094 * [Math]::Pow(10, 19)
$ws.Range("G2").Value = 1.354078223128095 * [Math]::Pow(10, 19)

# Row 3
$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 3.536033448013082
